$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Mayo de 2020 a las 21:05"
$ws.Cells.Item(4, 2).Value = 1737429
$ws.Cells.Item(4, 3).Value = 12154
$ws.Cells.Item(4, 4).Value = 483092
$ws.Cells.Item(4, 5).Value = 1152841
$ws.Cells.Item(4, 7).Value = 924
$ws.Cells.Item(4, 8).Value = 101496
$ws.Cells.Item(10, 2).Value = 182913
$ws.Cells.Item(10, 3).Value = 191
$ws.Cells.Item(10, 4).Value = 66584
$ws.Cells.Item(10, 5).Value = 87733
$ws.Cells.Item(10, 7).Value = 66
$ws.Cells.Item(10, 8).Value = 28596
$ws.Cells.Item(11, 2).Value = 181770
$ws.Cells.Item(11, 3).Value = 482
$ws.Cells.Item(11, 5).Value = 10444
$ws.Cells.Item(11, 7).Value = 28
$ws.Cells.Item(11, 8).Value = 8526
$ws.Cells.Item(13, 2).Value = 158086
$ws.Cells.Item(13, 3).Value = 7293
$ws.Cells.Item(13, 5).Value = 85803
$ws.Cells.Item(40, 2).Value = 19666
$ws.Cells.Item(40, 3).Value = 910
$ws.Cells.Item(40, 4).Value = 5205
$ws.Cells.Item(40, 5).Value = 13645
$ws.Cells.Item(40, 7).Value = 19
$ws.Cells.Item(40, 8).Value = 816
$ws.Cells.Item(77, 2).Value = 3369
$ws.Cells.Item(77, 3).Value = 79
$ws.Cells.Item(77, 5).Value = 687
$ws.Cells.Item(102, 1).Value = "Sri Lanka"
$ws.Cells.Item(102, 2).Value = 1469
$ws.Cells.Item(102, 3).Value = 150
$ws.Cells.Item(102, 4).Value = 732
$ws.Cells.Item(102, 5).Value = 727
$ws.Cells.Item(102, 8).Value = 10
$ws.Cells.Item(103, 1).Value = "Maldivas"
$ws.Cells.Item(103, 2).Value = 1457
$ws.Cells.Item(103, 3).Value = 19
$ws.Cells.Item(103, 4).Value = 197
$ws.Cells.Item(103, 5).Value = 1255
$ws.Cells.Item(103, 8).Value = 5
$ws.Cells.Item(105, 2).Value = 1195
$ws.Cells.Item(105, 3).Value = 17
$ws.Cells.Item(105, 5).Value = 1146
$ws.Cells.Item(116, 2).Value = 984
$ws.Cells.Item(116, 3).Value = 28
$ws.Cells.Item(116, 4).Value = 639
$ws.Cells.Item(116, 5).Value = 335
$ws.Cells.Item(129, 1).Value = "Republica del Chad"
$ws.Cells.Item(129, 2).Value = 715
$ws.Cells.Item(129, 3).Value = 15
$ws.Cells.Item(129, 4).Value = 359
$ws.Cells.Item(129, 5).Value = 292
$ws.Cells.Item(129, 7).Value = 2
$ws.Cells.Item(129, 8).Value = 64
$ws.Cells.Item(130, 1).Value = "Crucero"
$ws.Cells.Item(130, 2).Value = 712
$ws.Cells.Item(130, 4).Value = 651
$ws.Cells.Item(130, 5).Value = 48
$ws.Cells.Item(130, 8).Value = 13
$ws.Cells.Item(167, 1).Value = "Zimbabue"
$ws.Cells.Item(167, 2).Value = 132
$ws.Cells.Item(167, 3).Value = 76
$ws.Cells.Item(167, 4).Value = 25
$ws.Cells.Item(167, 5).Value = 103
$ws.Cells.Item(167, 8).Value = 4
$ws.Cells.Item(168, 1).Value = "Camboya"
$ws.Cells.Item(168, 2).Value = 124
$ws.Cells.Item(168, 4).Value = 122
$ws.Cells.Item(168, 5).Value = 2
$ws.Cells.Item(168, 8).Value = 0
$ws.Cells.Item(169, 1).Value = "Siria"
$ws.Cells.Item(169, 2).Value = 121
$ws.Cells.Item(169, 4).Value = 43
$ws.Cells.Item(169, 5).Value = 74
$ws.Cells.Item(169, 8).Value = 4
$ws.Cells.Item(170, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(170, 2).Value = 116
$ws.Cells.Item(170, 4).Value = 108
$ws.Cells.Item(170, 5).Value = 0
$ws.Cells.Item(170, 8).Value = 8
$ws.Cells.Item(171, 1).Value = "Malaui"
$ws.Cells.Item(171, 4).Value = 37
$ws.Cells.Item(171, 5).Value = 60
$ws.Cells.Item(171, 8).Value = 4
$ws.Cells.Item(172, 1).Value = "Aruba"
$ws.Cells.Item(172, 2).Value = 101
$ws.Cells.Item(172, 4).Value = 97
$ws.Cells.Item(172, 5).Value = 1
$ws.Cells.Item(172, 8).Value = 3
$ws.Cells.Item(173, 1).Value = "Bahamas"
$ws.Cells.Item(173, 2).Value = 100
$ws.Cells.Item(173, 4).Value = 46
$ws.Cells.Item(173, 5).Value = 43
$ws.Cells.Item(173, 8).Value = 11
$ws.Cells.Item(174, 1).Value = "Monaco"
$ws.Cells.Item(174, 2).Value = 98
$ws.Cells.Item(174, 4).Value = 90
$ws.Cells.Item(174, 5).Value = 4
$ws.Cells.Item(174, 8).Value = 4
$ws.Cells.Item(175, 1).Value = "Barbados"
$ws.Cells.Item(175, 2).Value = 92
$ws.Cells.Item(175, 4).Value = 71
$ws.Cells.Item(175, 5).Value = 14
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 7
$ws.Cells.Item(176, 1).Value = "Comoras"
$ws.Cells.Item(176, 2).Value = 87
$ws.Cells.Item(176, 4).Value = 24
$ws.Cells.Item(176, 5).Value = 61
$ws.Cells.Item(176, 7).Value = 1
$ws.Cells.Item(176, 8).Value = 2
$ws.Cells.Item(177, 1).Value = "Liechtenstein"
$ws.Cells.Item(177, 2).Value = 82
$ws.Cells.Item(177, 3).Value = 0
$ws.Cells.Item(177, 4).Value = 55
$ws.Cells.Item(177, 5).Value = 26
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(178, 1).Value = "Libia"
$ws.Cells.Item(178, 2).Value = 79
$ws.Cells.Item(178, 3).Value = 2
$ws.Cells.Item(178, 4).Value = 40
$ws.Cells.Item(178, 5).Value = 36
$ws.Cells.Item(178, 8).Value = 3
$ws.Cells.Item(179, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(179, 2).Value = 77
$ws.Cells.Item(179, 3).Value = 0
$ws.Cells.Item(179, 4).Value = 60
$ws.Cells.Item(179, 5).Value = 2
$ws.Cells.Item(179, 8).Value = 15
$ws.Cells.Item(180, 1).Value = "Angola"
$ws.Cells.Item(180, 2).Value = 71
$ws.Cells.Item(180, 3).Value = 1
$ws.Cells.Item(180, 4).Value = 18
$ws.Cells.Item(180, 5).Value = 49
$ws.Cells.Item(180, 8).Value = 4
$ws.Cells.Item(181, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(181, 2).Value = 60
$ws.Cells.Item(181, 4).Value = 60
$ws.Cells.Item(181, 5).Value = 0
$ws.Cells.Item(181, 8).Value = 0
$ws.Cells.Item(199, 1).Value = "Belice"
$ws.Cells.Item(199, 4).Value = 16
$ws.Cells.Item(199, 8).Value = 2
$ws.Cells.Item(201, 1).Value = "Santa Lucia"
$ws.Cells.Item(201, 4).Value = 18
$ws.Cells.Item(201, 8).Value = 0
$ws.Cells.Item(213, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(213, 4).Value = 7
$ws.Cells.Item(213, 8).Value = 1
$ws.Cells.Item(214, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 4).Value = 8
$ws.Cells.Item(214, 8).Value = 0
